$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 796 (shifts existing rows 796-837 down to 797-838),
# then populate it with the new entry: 2026/02/15, 日, 4, 22
$ws.Rows.Item(796).Insert()

# Force column A to be treated as text so the date-looking string isn't
# auto-converted into a date serial number (matches the existing rows,
# which store dates as literal text).
$ws.Cells.Item(796, 1).NumberFormat = "@"
$ws.Cells.Item(796, 1).Value = "2026/02/15"
$ws.Cells.Item(796, 1).ClearFormats()

$ws.Cells.Item(796, 2).Value = "日"
$ws.Cells.Item(796, 3).Value = 4
$ws.Cells.Item(796, 4).Value = 22
